# Apply the diff: insert a new data row at sheet row 636 (shifting existing
# rows 636..726 down to 637..727), and populate the new row with the data
# that was recorded for that date. The carried-over columns (A,B,C,E,F,G,H,
# I,J,K,L,Q,T) keep the same values as the row that used to occupy position
# 636 (now at 637), while D,M,N,O,P,R,S get the new values for the inserted
# record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 636, pushing old rows 636-726 to 637-727.
$ws.Rows.Item(636).Insert()

# Fill in the carried-over (unchanged) columns for the new row 636, copying
# the values that remain constant for this data template.
$ws.Cells.Item(636, 1).Value  = 9                                          # A: Mercado ID
$ws.Cells.Item(636, 2).Value  = "Vega Central Mapocho de Santiago"          # B: Mercado
$ws.Cells.Item(636, 3).Value  = "Metropolitana"                             # C: Región
$ws.Cells.Item(636, 5).Value  = 13                                          # E: Codreg
$ws.Cells.Item(636, 6).Value  = "Fruta"                                     # F: Tipo
$ws.Cells.Item(636, 7).Value  = 100108                                      # G: Producto ID
$ws.Cells.Item(636, 8).Value  = "Tropicales y subtropicales"                # H: Producto
$ws.Cells.Item(636, 9).Value  = 100108002                                   # I: Categoría ID
$ws.Cells.Item(636, 10).Value = "Mango"                                     # J: Categoría
$ws.Cells.Item(636, 11).Value = "Sin especificar"                           # K: Variedad
$ws.Cells.Item(636, 12).Value = "Primera"                                   # L: Calidad
$ws.Cells.Item(636, 17).Value = "$/bandeja 4 kilos"                         # Q: Unidad de comercialización
$ws.Cells.Item(636, 20).Value = 4                                           # T: Kg / unidad

# New values specific to the inserted record.
$ws.Cells.Item(636, 4).Value2 = 45154    # D: Fecha
$ws.Cells.Item(636, 4).NumberFormat = $ws.Cells.Item(637, 4).NumberFormat
$ws.Cells.Item(636, 13).Value = 520      # M: Volumen
$ws.Cells.Item(636, 14).Value = 8000     # N: Precio mínimo
$ws.Cells.Item(636, 15).Value = 8500     # O: Precio máximo
$ws.Cells.Item(636, 16).Value = 8212     # P: Precio promedio ponderado
$ws.Cells.Item(636, 18).Value = "Brasil" # R: Origen
$ws.Cells.Item(636, 19).Value = 2053     # S: Precio $/Kg
